$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting old rows 145-146 down to 146-147.
$ws.Rows.Item(145).Insert()

# Populate the new row 145 with the new weekly price record.
$ws.Range("A145").Value = 4
$ws.Range("B145").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C145").Value = "Los Lagos"
$ws.Range("D145").Value = 44448
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 100112044
$ws.Range("G145").Value = "Perejil"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 80
$ws.Range("K145").Value = 5000
$ws.Range("L145").Value = 5000
$ws.Range("M145").Value = 5000
$ws.Range("N145").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O145").Value = "Región Metropolitana"
$ws.Range("P145").Value = 1667
$ws.Range("Q145").Value = 3
$ws.Range("R145").Value = "Hortaliza"
